$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = '27.480.28'
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +3.31%  '
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = '1.820.33'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +4.98%  '
$ws.Cells.Item(4, 5).Value = '  +0.69%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '343.65'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +2.83%  '
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '1.000'
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.65%  '
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = '0.3844'
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +3.46%  '
$ws.Cells.Item(8, 5).Value = '  +4.29%  '
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = '48.95'
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -1.28%  '
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '1.237'
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +2.75%  '
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '0.07794'
$c.Style = "Normal"
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.78%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = '22.36'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +9.12%  '
$ws.Cells.Item(14, 5).Value = '  +2.86%  '
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = '1.819.04'
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +5.35%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = '7.228'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +2.99%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = '0.00001123'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +2.59%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value = '0.06719'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +0.60%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = '86.40'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +3.37%  '
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.61%  '
$ws.Cells.Item(21, 5).Value = '  +5.28%  '
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '6.577'
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +6.48%  '
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = '13.23'
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.45%  '
$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '27.472.64'
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +3.50%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = '2.463'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.03%  '
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '2.694'
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +6.80%  '
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = '22.36'
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +14.10%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '1.472'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +3.17%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '154.05'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.34%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = '2.021.14'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +5.27%  '
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = '136.53'
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +3.70%  '
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '6.382'
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +3.83%  '
$ws.Cells.Item(33, 5).Value = '  -1.17%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '13.93'
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +5.50%  '
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = '0.08814'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +2.56%  '
$ws.Cells.Item(36, 5).Value = '  -0.80%  '
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = '5.630'
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +3.23%  '
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '0.7043'
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +12.48%  '
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = '0.2266'
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +4.26%  '
$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '0.02410'
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +2.35%  '
$ws.Cells.Item(41, 2).Value = 'Hedera'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = '0.06496'
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +2.15%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '8.988'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.60%  '
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '1.298'
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +4.45%  '
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = '14.78'
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.64%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '0.6598'
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +8.63%  '
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.9997'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +0.57%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '3.960'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +1.50%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = '2.192'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +6.18%  '
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '132.75'
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +2.46%  '
$ws.Cells.Item(50, 5).Value = '  -0.04%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '80.80'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +3.63%  '
